$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 ("Save"), matching the style used by the other header cells (e.g. G1)
$ws.Cells.Item(1, 8).Value = "Save"
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)  # xlPasteFormats

# Data cells H2:H8 = 0 (plain numeric, default style)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
